$d = $word.ActiveDocument

# The very first paragraph holds the AFFARS topic-id placeholder, spread
# across two runs: "**ID__AFFARS_mp_5315_3_topic_34__ID**" + a trailing
# space run. Replace that whole stretch with the new placeholder text in
# a single Find/Replace so the now-unneeded trailing-space run disappears.
$p1 = $d.Paragraphs(1)
$p1.Range.Find.Execute( `
    "**ID__AFFARS_mp_5315_3_topic_34__ID** ", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "**ID__AFFARS_MP_5315_3_6_4__ID**", 2)

# Re-fetch the paragraph (defensive - Find/Replace can reshuffle ranges)
# and give it a (borderless) paragraph border, i.e. just padding/space of
# 5 on every side, plus bump the left indent from 120 to 225 twips
# (6pt -> 11.25pt).
$p1 = $d.Paragraphs(1)
$p1.Format.Borders.DistanceFromTop = 5
$p1.Format.Borders.DistanceFromLeft = 5
$p1.Format.Borders.DistanceFromBottom = 5
$p1.Format.Borders.DistanceFromRight = 5
$p1.Format.LeftIndent = 11.25
